# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new labels in AD1:AF1, matching the existing
#     header formatting (bold font, thin border, centered/top aligned). ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerStyle = $ws.Range("AC1:AF1")
$headerStyle.Font.Bold = $true
$headerStyle.HorizontalAlignment = -4108
$headerStyle.VerticalAlignment = -4160
$headerStyle.Borders.LineStyle = 1
$headerStyle.Borders.Weight = 2

# --- Data rows (2-56): every player gets the same 2022 Yankees season
#     record -- 99 wins, 63 losses, 0 ties. ---
$lastRow = 56
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 99
    $ws.Cells.Item($r, 31).Value = 63
    $ws.Cells.Item($r, 32).Value = 0
}
